# Update Unit of Measure values in row 6 of the YROH sheet from
# "G" / "PC" to "KG", and align the formatting of the previously
# left-aligned UoM cells (UoM - Primary/Secondary/Pricing/Production/
# Purchasing/Shipping/Component) with the style already used by the
# neighboring "Unit of Weight"/"Base UoM" cells (K6/N6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YROH")

# Copy K6's format (no explicit horizontal alignment) onto the cells
# whose style needs to drop the left-alignment override, reusing the
# existing cell style instead of creating a new one.
$ws.Range("K6").Copy()
$styleCols = @("W", "X", "Y", "Z", "AE", "AF", "AG")
foreach ($col in $styleCols) {
    $ws.Range("$col`6").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Update the cell values to "KG" for all affected Unit of Measure columns.
$cols = @("K", "N", "W", "X", "Y", "Z", "AE", "AF", "AG")
foreach ($col in $cols) {
    $ws.Range("$col`6").Value = "KG"
}

# Move the active selection to A6, matching the saved view state.
$ws.Range("A6").Select()
